$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Replace row 6 ("Jumper2" / JP1,JP2 / PIN_HEADER_1x2) with the new
#    "Jumper3" single 3-pin header (JP1) component.
#    We stage the new text in scratch cells, then Copy + PasteSpecial
#    (values only) into the destination cells so the destination keeps its
#    existing cell formatting (border + quote-prefix) and only the value
#    changes - this mirrors how the text was edited in the original
#    worksheet (formatting preserved, just the designator swapped out).
#    Scratch cells are written in the same order the new strings should be
#    appended to the shared-string table: Jumper3, PIN_HEADER_1x3, JP1.
# ---------------------------------------------------------------------------
$ws.Range("Z1").Value = "Jumper3"
$ws.Range("Z2").Value = "PIN_HEADER_1x3"
$ws.Range("Z3").Value = "JP1"

$ws.Range("Z1").Copy()
$ws.Range("A6").PasteSpecial(-4163)

$ws.Range("Z2").Copy()
$ws.Range("D6").PasteSpecial(-4163)

$ws.Range("Z3").Copy()
$ws.Range("C6").PasteSpecial(-4163)

$ws.Range("Z1").Copy()
$ws.Range("E6").PasteSpecial(-4163)

$ws.Range("Z1:Z3").Clear()

# Numeric columns for the new row: Qty 1, Unit Price 0.13, Total Price 1.3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.13
$ws.Range("H6").Value = 1.3

# ---------------------------------------------------------------------------
# 2. Colour-code the "Comment" column (A) to flag each component category.
#    Accent1 (blue)  -> capacitor / connectors / encoder rows (2-5)
#    Accent5 (lt blue)-> newly added jumper + first-of-group resistor rows
#    Red             -> remaining resistor rows / the IC row
# ---------------------------------------------------------------------------
$ws.Range("A2").Font.ThemeColor = 5
$ws.Range("A3").Font.ThemeColor = 5
$ws.Range("A4").Font.ThemeColor = 5
$ws.Range("A5").Font.ThemeColor = 5

$ws.Range("A6").Font.ThemeColor = 9
$ws.Range("A7").Font.ThemeColor = 9
$ws.Range("A9").Font.ThemeColor = 9

$ws.Range("A8").Font.Color = 255
$ws.Range("A10").Font.Color = 255

# ---------------------------------------------------------------------------
# 3. Column widths: give column B its own best-fit width, keep A and C:F at
#    the original 20.140625 width (previously A:F shared one column group).
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 48

# ---------------------------------------------------------------------------
# 4. Selection cursor moved to C31 (reflects the last cell the author had
#    selected when saving).
# ---------------------------------------------------------------------------
$ws.Range("C31").Select()
